# F05 Froze Encoder 12345
# Update the "Epoch Accuracy" sheet with refreshed accuracy values coming
# from a re-run of the training notebook (new DisplayOutputs object
# identity + updated accuracy numbers in column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: updated accuracy values ---
$ws.Range("B2").Value = 0.453125
$ws.Range("B3").Value = 0.3125
$ws.Range("B5").Value = 0.234375
$ws.Range("B6").Value = 0.25
$ws.Range("B8").Value = 0.21875
$ws.Range("B9").Value = 0.203125
$ws.Range("B10").Value = 0.234375
$ws.Range("B11").Value = 0.25
$ws.Range("B12").Value = 0.234375
$ws.Range("B14").Value = 0.171875
$ws.Range("B15").Value = 0.171875
$ws.Range("B16").Value = 0.171875
$ws.Range("B17").Value = 0.15625
$ws.Range("B18").Value = 0.15625
$ws.Range("B19").Value = 0.15625
$ws.Range("B20").Value = 0.15625
$ws.Range("B22").Value = 0.15625
$ws.Range("B23").Value = 0.15625
$ws.Range("B24").Value = 0.15625
$ws.Range("B25").Value = 0.15625
$ws.Range("B26").Value = 0.15625
$ws.Range("B27").Value = 0.140625
$ws.Range("B28").Value = 0.140625
$ws.Range("B29").Value = 0.140625
$ws.Range("B30").Value = 0.140625
$ws.Range("B31").Value = 0.140625
$ws.Range("B32").Value = 0.140625
$ws.Range("B33").Value = 0.140625
$ws.Range("B34").Value = 0.140625
$ws.Range("B35").Value = 0.140625
$ws.Range("B36").Value = 0.140625
$ws.Range("B37").Value = 0.140625
$ws.Range("B38").Value = 0.140625
$ws.Range("B39").Value = 0.140625
$ws.Range("B44").Value = 0.140625
$ws.Range("B45").Value = 0.140625
$ws.Range("B46").Value = 0.140625
$ws.Range("B47").Value = 0.140625
$ws.Range("B48").Value = 0.140625
$ws.Range("B103").Value = 0.109375
$ws.Range("B104").Value = 0.109375
$ws.Range("B105").Value = 0.140625
$ws.Range("B106").Value = 0.125
$ws.Range("B107").Value = 0.0625
$ws.Range("B109").Value = 0.078125
$ws.Range("B112").Value = 0.109375
$ws.Range("B113").Value = 0.09375
$ws.Range("B114").Value = 0.109375
$ws.Range("B115").Value = 0.125
$ws.Range("B116").Value = 0.125
$ws.Range("B117").Value = 0.0625
$ws.Range("B118").Value = 0.08196721311475409

# --- Column A (rows 102-118): repr() of the re-created DisplayOutputs
#     object, whose memory address changed because the notebook cell
#     was re-executed. ---
$newRepr = "<__main__.DisplayOutputs object at 0x7f71a0717b50>"
for ($r = 102; $r -le 118; $r++) {
    $ws.Cells.Item($r, 1).Value = $newRepr
}
